# Fixing flags in school districts
# The "Changed Values" sheet originally listed 25 changed ids (rows 2-26).
# After the fix, only the 7 rows that correspond to school-district flag
# changes (the former rows 20-26) remain. For those rows the x1_old/x1_new
# values (columns B/C) were swapped, and the rows were re-sorted by the
# (corrected) x1_old value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changed Values")

# Remove the obsolete rows (old ids rows 2-19), shifting rows 20-26 up to 2-8.
$ws.Rows("2:19").Delete()

# Final corrected values for the remaining 7 rows (id, x1_old, x1_new),
# already swapped and sorted by the new x1_old value.
$finalData = @(
    @(40620,   20, 22),
    @(33447,   21, 23),
    @(1265711, 22, 24),
    @(1036991, 23, 20),
    @(38235,   24, 25),
    @(38453,   25, 26),
    @(64951,   26, 21)
)

for ($i = 0; $i -lt $finalData.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $finalData[$i][0]
    $ws.Cells.Item($row, 2).Value = $finalData[$i][1]
    $ws.Cells.Item($row, 3).Value = $finalData[$i][2]
    $ws.Cells.Item($row, 4).Value = "x1"
}
